$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = "Update Only"
$ws.Range("R2").Value = "No"
$ws.Range("R3").Value = "No"
$ws.Range("R4").Value = "No"
$ws.Range("R5").Value = "No"
$ws.Range("R6").Value = "No"
$ws.Range("R7").Value = "No"
$ws.Range("R8").Value = "No"
$ws.Range("R9").Value = "No"

$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("R3:R9").Select()
